$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "국내에서 유일하게 제대로 Data Science 를 가르치는 학부"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/snu-computational-science-faculty/#utm_source=rss&utm_medium=rss&utm_campaign=snu-computational-science-faculty"

$ws.Range("D16").Value = "7. Time series 분석 정리 (1)"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/164"

$ws.Range("D50").Value = "볼츠만 머신 [설명]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7520379"
